$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B width: 16.42578125 -> 32 ---
# ColumnWidth setter applies an internal +0.8333 pixel-rounding offset in this
# engine, so request a value that lands exactly on 32 after that adjustment.
$ws.Columns("B").ColumnWidth = 31.166666666666668

# --- Row 1: C1 becomes bold "Condiciones:" (same look as header row) ---
$ws.Range("C1").Value = "Condiciones:"
$ws.Range("C1").Font.Bold = $true

# --- Row 4: add explanatory C4 (written before C2/C3 to match shared-string order) ---
$ws.Range("C4").Value = "Para los filtros y cálculos se utiliza la fecha de envío de las tareas. Solo se consideran días laborales entre 9 y 18 horas."

# --- Row 2: keep A2 text, add empty left-aligned B2, add explanatory C2 ---
$ws.Range("B2").HorizontalAlignment = -4131  # xlLeft (matches A2/A3/A4 label style)
$ws.Range("C2").Value = "Se consideran procesos de gestión documental generados fuera de la unidad de permanencia y con estado no anulado."

# --- Row 3: add explanatory C3 ---
$ws.Range("C3").Value = "Se consideran procesos que fueron atendidos por la unidad de permanencia. No se consideran procesos cerrados en la unidad de permanencia."

# --- Date cells B3:B4 gain left horizontal alignment while keeping their
#     original m/d/yyyy (numFmtId 14) date format ---
$dateRange = $ws.Range("B3")
$dateRange.HorizontalAlignment = -4131  # xlLeft
$dateRange.NumberFormat = "mm-dd-yy"    # re-assert builtin short-date format (id 14)
$dateRange.Copy() | Out-Null
$ws.Range("B4").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

$wb.Save()
